# Commit: "Ich habe GitHub geputzt." — fix the typo'd package name so it
# reads "Complex \ Heatmap" instead of "ComplexHeatmap" in the packages table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A11").Value = "Complex \ Heatmap"

# Reflect the author's final cursor/scroll position in the saved view state.
$ws.Activate()
$ws.Range("A11").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
